# Converted shape files from growers
# Adds a "DryMoisture_perc" column (C) to the BulkDensity sheet, updates a
# couple of Bu_t conversion values, and appends a new "Blow Out" row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("BulkDensity")

# --- New "Blow Out" row (added first so shared-string ordering matches) ----
$ws.Range("A17").Value = "Blow Out"
$ws.Range("A17").Style = $ws.Range("A16").Style
$ws.Range("B17").Value = "NA"
$ws.Range("C17").Value = "NA"

# --- Header ---------------------------------------------------------------
$ws.Range("C1").Value = "DryMoisture_perc"

# --- Corrected Bu_t values (Canola / Mustard rows) -------------------------
$ws.Range("B3").Value = 42.4
$ws.Range("B7").Value = 42.4
$ws.Range("B14").Value = 38.348513513513517

# --- New DryMoisture_perc values -------------------------------------------
$ws.Range("C2").Value = 10
$ws.Range("C3").Value = 6.6
$ws.Range("C4").Value = 10
$ws.Range("C5").Value = 9.7
$ws.Range("C6").Value = 9.8
$ws.Range("C6").WrapText = $true
$ws.Range("C7").Value = 6.6
$ws.Range("C8").ClearFormats()
$ws.Range("C8").Value = 13
$ws.Range("C9").Value = 10
$ws.Range("C10").Value = 7
$ws.Range("C11").Value = 7
$ws.Range("C11").WrapText = $true
$ws.Range("C12").Value = 9.8
$ws.Range("C12").WrapText = $true
$ws.Range("C13").Value = 9.8
$ws.Range("C13").WrapText = $true
$ws.Range("C14").ClearFormats()
$ws.Range("C14").Value = 14
$ws.Range("C15").Value = 13
$ws.Range("C15").WrapText = $true
$ws.Range("C16").Value = 12.5
$ws.Range("C16").WrapText = $true

# --- Move the active selection to match the author's final cursor position -
$ws.Range("B15").Select()
